$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New log entries for 9 March 2021 (serial date 44264), appended as rows
# 14-17 beneath the existing data (which currently ends at row 12).
# ---------------------------------------------------------------------------

# Row 14: Making scene, adding rigidbody, checking specfic boxes are ticked
$ws.Range("A14").Value = "Making scene, adding rigidbody, checking specfic boxes are ticked"
$ws.Range("B14").Value = 44264
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C14").Value = 0.2
$ws.Range("D14").Value = 0.54861111111111105
$ws.Range("D2").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = 0.55902777777777779
$ws.Range("E2").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("G14").Value = 0.15

# Row 15: Adding player movement
$ws.Range("A15").Value = "Adding player movement "
$ws.Range("B15").Value = 44264
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").Value = 0.05
$ws.Range("D15").Value = 0.55902777777777779
$ws.Range("D2").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0.5625
$ws.Range("E2").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 0.05

# Row 16: Writing triggerRespawn script
$ws.Range("A16").Value = "Writing triggerRespawn script "
$ws.Range("B16").Value = 44264
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").Value = 0.3
$ws.Range("D16").Value = 0.5625
$ws.Range("D2").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0.59027777777777779
$ws.Range("E2").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 0.1
$ws.Range("G16").Value = 0.4

# Row 17: Adding triggerRespawn script
$ws.Range("A17").Value = "Adding triggerRespawn script "
$ws.Range("B17").Value = 44264
$ws.Range("B2").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C17").Value = 0.05
$ws.Range("D17").Value = 0.59027777777777779
$ws.Range("D2").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = 0.59791666666666665
$ws.Range("E2").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 0.06
$ws.Range("G17").Value = 0.11

# ---------------------------------------------------------------------------
# Column A got wider to fit the new, longer task descriptions.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 58.5

# ---------------------------------------------------------------------------
# Move the active-cell selection down below the newly added data, as Excel
# would leave it after the user finished typing in the sheet.
# ---------------------------------------------------------------------------
[void]$ws.Range("A18").Select()
